$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 11
$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "React Developer"
$ws.Cells.Item($row, 3).Value = "dgfsegse"
$ws.Cells.Item($row, 4).Value = 2
$ws.Cells.Item($row, 5).Value = 5
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
